$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Build a "template" cell far off the used range that carries exactly
# the formatting we want the H column to end up with: Text number format
# ("@"), Arial 10pt black font, holding the literal text "False" (not the
# boolean). A leading apostrophe forces Excel to store it as text even
# though it reads like a boolean keyword.
$tpl = $ws.Cells.Item(200, 200)
$tpl.NumberFormat = "@"
$tpl.Font.Name = "Arial"
$tpl.Font.Size = 10
$tpl.Font.Color = 0
$tpl.Value = "'False"

# Copy that formatted text value onto every row of column H that currently
# holds the boolean FALSE (rows 2-27), replacing FALSE -> "False" text.
for ($r = 2; $r -le 27; $r++) {
    $tpl.Copy()
    $ws.Cells.Item($r, 8).PasteSpecial(-4163)
}

# Remove the scratch/template cell entirely so it doesn't linger in the
# used range / dimension.
$tpl.Clear()

# Re-apply the number format + font once more directly on the real H
# range so every cell (including ones that started out with a different
# style, e.g. H2) ends up on the identical style record.
$rng = $ws.Range("H2:H32")
$rng.NumberFormat = "@"
$rng.Font.Name = "Arial"
$rng.Font.Size = 10
$rng.Font.Color = 0

$ws.Application.CutCopyMode = $false

# Match the saved view/selection state.
$ws.Range("A25").Select()
$ws.Range("G28:J34").Select()
